$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the underlying data values that feed the two chart series
# (Branch and Bounding -> column C, Greedy 2-approximation -> column D)
$ws.Range("C4").Value = 23.722300000000001
$ws.Range("D4").Value = 26.508500000000002

$ws.Range("C5").Value = 200.6523
$ws.Range("D5").Value = 225.26849999999999

$ws.Range("C6").Value = 157.27529999999999
$ws.Range("D6").Value = 170.8321

$ws.Range("C7").Value = 233.66
$ws.Range("D7").Value = 267.31

$ws.Range("C8").Value = 371.09640000000002
$ws.Range("D8").Value = 444.2998

$ws.Range("C9").Value = 376.25450000000001
$ws.Range("D9").Value = 418.84820000000002

$ws.Range("C10").Value = 387.69389999999999
$ws.Range("D10").Value = 405.78190000000001

$ws.Range("C11").Value = 394.65120000000002
$ws.Range("D11").Value = 466.03190000000001

$ws.Range("C12").Value = 113.40179999999999
$ws.Range("D12").Value = 129.03370000000001

$ws.Range("C13").Value = 140.398
$ws.Range("D13").Value = 155.43559999999999

$ws.Range("C14").Value = 141.2585
$ws.Range("D14").Value = 159.94800000000001

# Update the selection shown on the sheet
$ws.Range("C4:D14").Select()
